$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.852.55'
$ws.Range("E2").Value = '  +0.06%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.535.89'
$ws.Range("E3").Value = '  +0.28%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.08'
$ws.Range("E5").Value = '  -1.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '195.68'
$ws.Range("E6").Value = '  +4.67%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.205'
$ws.Range("E9").Value = '  -4.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.651'
$ws.Range("E10").Value = '  -0.98%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.75'
$ws.Range("E11").Value = '  +0.53%  '

$ws.Range("E12").Value = '  -1.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.54'
$ws.Range("E13").Value = '  -1.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.101.34'
$ws.Range("E14").Value = '  +0.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '605.95'
$ws.Range("E15").Value = '  -1.56%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '70.059.54'
$ws.Range("E16").Value = '  +0.19%  '

$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.16'
$ws.Range("E17").Value = '  +0.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.70'
$ws.Range("E18").Value = '  -0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.539.94'
$ws.Range("E19").Value = '  -0.63%  '

$ws.Range("E20").Value = '  +0.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.993'
$ws.Range("E21").Value = '  -0.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.17'
$ws.Range("E22").Value = '  +3.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.22'
$ws.Range("E23").Value = '  +3.61%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.60'
$ws.Range("E24").Value = '  -2.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.61'
$ws.Range("E25").Value = '  -1.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.15'
$ws.Range("E26").Value = '  +4.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.91'
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.63'
$ws.Range("E28").Value = '  -3.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.53'
$ws.Range("E29").Value = '  -1.53%  '

$ws.Range("B30").Value = 'dogwifhat'
$ws.Range("C30").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.34'
$ws.Range("E30").Value = '  +16.29%  '

$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.12'
$ws.Range("E31").Value = '  +1.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.59'
$ws.Range("E32").Value = '  +1.35%  '

$ws.Range("E33").Value = '  -1.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.15'
$ws.Range("E34").Value = '  -0.88%  '

$ws.Range("D35").Value = '0.0₃0859'
$ws.Range("E35").Value = '  +10.83%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.734.24'
$ws.Range("E36").Value = '  +5.18%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'

$ws.Range("E38").Value = '  -2.64%  '

$ws.Range("E39").Value = '  +1.48%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.393'
$ws.Range("E40").Value = '  -1.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.59'
$ws.Range("E41").Value = '  -0.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '488.62'
$ws.Range("E42").Value = '  -8.76%  '

$ws.Range("E43").Value = '  -5.29%  '

$ws.Range("E44").Value = '  -0.32%  '

$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.84'
$ws.Range("E45").Value = '  -3.81%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.141'
$ws.Range("E46").Value = '  -1.96%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.30'
$ws.Range("E47").Value = '  -1.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.59'
$ws.Range("E49").Value = '  -4.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000254'
$ws.Range("E50").Value = '  +5.90%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.65'
$ws.Range("E51").Value = '  -1.20%  '
